$d = $word.ActiveDocument

# 1. Fix the typo "retrive" -> "retrieve" in the "Users can add, ..." sentence.
$d.Content.Find.Execute("retrive", $true, $false, $false, $false, $false,
                         $true, 1, $false, "retrieve", 2)

# 2. Remove the old "_GoBack" bookmark (between "desk" and "top with C#.net").
if ($d.Bookmarks.Exists("_GoBack")) {
    $d.Bookmarks.Item("_GoBack").Delete()
}

# 3. Re-add the "_GoBack" bookmark at the end of the sentence we just edited
#    ("Users can add, retrieve and delete their messages. ").
$findRange = $d.Content
$findRange.Find.Execute("Users can add, retrieve and delete their messages. ",
                         $true, $false, $false, $false, $false,
                         $true, 1, $false, "", 0)
$gobackRange = $findRange.Duplicate
$gobackRange.Collapse(0)
$d.Bookmarks.Add("_GoBack", $gobackRange)
